$wb = $excel.ActiveWorkbook

# --- Sheet "Errores": insert two new rows (new findings) and two new shared strings ---
$wsErrores = $wb.Worksheets.Item("Errores")

# Insert two blank rows above the current row 13 (AuctionManagementBean block),
# pushing it down to row 15 and the BussinessException block down to row 17.
$wsErrores.Rows.Item(13).EntireRow.Insert() | Out-Null
$wsErrores.Rows.Item(13).EntireRow.Insert() | Out-Null

# Fill the two newly-inserted rows with the new comments (column C only).
$wsErrores.Cells.Item(12, 3).Value = "Muchos de los metodos no usan la implementacion del toBO por tanto se pueden generar muchos errores al momento de mapear la info"
$wsErrores.Cells.Item(13, 3).Value = "El metodo de cerrar subasta no esta persistiendo los cambios"

# Page setup for the "Errores" sheet.
$wsErrores.PageSetup.PaperSize = 9
$wsErrores.PageSetup.Orientation = 1

# Update the selection on the "Errores" sheet.
$wsErrores.Range("C11").Select() | Out-Null

# --- Sheet "Datos": becomes the active sheet/tab, with a new selection ---
$wsDatos = $wb.Worksheets.Item("Datos")
$wsDatos.Activate() | Out-Null
$wsDatos.Range("A38").Select() | Out-Null
